$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.7234583817159717
$ws.Range("D2").Value = 0.02992650373184347
$ws.Range("E2").Value = 0.04517576885281915
$ws.Range("F2").Value = 1.22091545983676
$ws.Range("G2").Value = 0.9908292259809173
$ws.Range("H2").Value = 0.01559863828661262
$ws.Range("I2").Value = 0.03268400461608323
$ws.Range("J2").Value = 0.6855829985444899
$ws.Range("K2").Value = 0.7033003184241338
$ws.Range("L2").Value = 0.08060822389467859
$ws.Range("M2").Value = 1.504576312338457
$ws.Range("N2").Value = 0.371838670746456
$ws.Range("B3").Value = 0.6322501639469067
$ws.Range("D3").Value = 0.02729660353406516
$ws.Range("E3").Value = 0.04228344199550405
$ws.Range("F3").Value = 1.166693831923979
$ws.Range("G3").Value = 0.9399337207665752
$ws.Range("H3").Value = 0.01952670636122256
$ws.Range("I3").Value = 0.03725143564904121
$ws.Range("J3").Value = 0.6672430924445081
$ws.Range("K3").Value = 0.6977634867156404
$ws.Range("L3").Value = 0.07560554798662222
$ws.Range("M3").Value = 1.308690821200059
$ws.Range("N3").Value = 0.3251636496990642
$ws.Range("B4").Value = 0.5759427315741732
$ws.Range("D4").Value = 0.02569338164078161
$ws.Range("E4").Value = 0.04049975012534324
$ws.Range("F4").Value = 1.133985308724967
$ws.Range("G4").Value = 0.9092956163280661
$ws.Range("H4").Value = 0.02225170067678983
$ws.Range("I4").Value = 0.04034215387394946
$ws.Range("J4").Value = 0.6563804476457165
$ws.Range("K4").Value = 0.6943434739125784
$ws.Range("L4").Value = 0.07251450665596249
$ws.Range("M4").Value = 1.189095425295932
$ws.Range("N4").Value = 0.2968317104157023
$ws.Range("B5").Value = 0.5522332016512621
$ws.Range("D5").Value = 0.02504712173110235
$ws.Range("E5").Value = 0.03974373909671902
$ws.Range("F5").Value = 1.119488908044346
$ws.Range("G5").Value = 0.8956971623432537
$ws.Range("H5").Value = 0.02344493764299449
$ws.Range("I5").Value = 0.04177919735687619
$ws.Range("J5").Value = 0.6514005433033958
$ws.Range("K5").Value = 0.6919363725644416
$ws.Range("L5").Value = 0.071225497212005
$ws.Range("M5").Value = 1.141674636252617
$ws.Range("N5").Value = 0.2859478655326058
$ws.Range("B6").Value = 0.5474528526684992
$ws.Range("D6").Value = 0.02494593903823272
$ws.Range("E6").Value = 0.03958634949663065
$ws.Range("F6").Value = 1.115495558945085
$ws.Range("G6").Value = 0.8919132737732411
$ws.Range("H6").Value = 0.02365598634535226
$ws.Range("I6").Value = 0.04216057768127257
$ws.Range("J6").Value = 0.6497869979392306
$ws.Range("K6").Value = 0.690307108562223
$ws.Range("L6").Value = 0.07098405708350208
$ws.Range("M6").Value = 1.135237953781854
$ws.Range("N6").Value = 0.2848633675745162
$ws.Range("B7").Value = 0.5733248250697898
$ws.Range("D7").Value = 0.02570112935661228
$ws.Range("E7").Value = 0.04040318447547175
$ws.Range("F7").Value = 1.129438008589517
$ws.Range("G7").Value = 0.9049233905299019
$ws.Range("H7").Value = 0.02228926641339557
$ws.Range("I7").Value = 0.04072104812072119
$ws.Range("J7").Value = 0.654152511177756
$ws.Range("K7").Value = 0.6909572038085017
$ws.Range("L7").Value = 0.07242367859553767
$ws.Range("M7").Value = 1.192366005279467
$ws.Range("N7").Value = 0.2986494062362084
$ws.Range("B8").Value = 0.6889996563998579
$ws.Range("D8").Value = 0.02903881049911305
$ws.Range("E8").Value = 0.04406581399966347
$ws.Range("F8").Value = 1.196353979809942
$ws.Range("G8").Value = 0.9676193037663978
$ws.Range("H8").Value = 0.01690984772249843
$ws.Range("I8").Value = 0.03462996388349726
$ws.Range("J8").Value = 0.6763323730585569
$ws.Range("K8").Value = 0.6969880918876221
$ws.Range("L8").Value = 0.07879183132047496
$ws.Range("M8").Value = 1.442092236922178
$ws.Range("N8").Value = 0.3582905455461542
$ws.Range("B9").Value = 0.9176783513777025
$ws.Range("D9").Value = 0.03568297661447417
$ws.Range("E9").Value = 0.05139662343027762
$ws.Range("F9").Value = 1.340444851673823
$ws.Range("G9").Value = 1.103464006143795
$ws.Range("H9").Value = 0.00895133212848756
$ws.Range("I9").Value = 0.02450132775710046
$ws.Range("J9").Value = 0.727073709974718
$ws.Range("K9").Value = 0.7146303827811664
$ws.Range("L9").Value = 0.09134639307886339
$ws.Range("M9").Value = 1.931189632465134
$ws.Range("N9").Value = 0.4746014894057566
$ws.Range("B10").Value = 1.079325035196604
$ws.Range("D10").Value = 0.04016979960233513
$ws.Range("E10").Value = 0.05614529199724361
$ws.Range("F10").Value = 1.42449394197223
$ws.Range("G10").Value = 1.184313372155003
$ws.Range("H10").Value = 0.005256224076847449
$ws.Range("I10").Value = 0.01893532468258119
$ws.Range("J10").Value = 0.7551175186008408
$ws.Range("K10").Value = 0.7137260700625916
$ws.Range("L10").Value = 0.1004910424845633
$ws.Range("M10").Value = 2.304118434087997
$ws.Range("N10").Value = 0.5526960559782168
$ws.Range("B11").Value = 1.114113439645791
$ws.Range("D11").Value = 0.0378175076925551
$ws.Range("E11").Value = 0.05537785145490215
$ws.Range("F11").Value = 1.259486817438912
$ws.Range("G11").Value = 1.038921543303047
$ws.Range("H11").Value = 0.02357112368085268
$ws.Range("I11").Value = 0.01827590859784944
$ws.Range("J11").Value = 0.6756998967127004
$ws.Range("K11").Value = 0.6107661674627494
$ws.Range("L11").Value = 0.1111774893655095
$ws.Range("M11").Value = 2.536825589285968
$ws.Range("N11").Value = 0.5020149088933721
$ws.Range("B12").Value = 1.112438246739657
$ws.Range("D12").Value = 0.03511121842335285
$ws.Range("E12").Value = 0.05777057773909089
$ws.Range("F12").Value = 1.116712561661345
$ws.Range("G12").Value = 0.9123595909832574
$ws.Range("H12").Value = 0.0624052142216982
$ws.Range("I12").Value = 0.01816197036427258
$ws.Range("J12").Value = 0.6091517994430404
$ws.Range("K12").Value = 0.5389493007312041
$ws.Range("L12").Value = 0.1284324305530546
$ws.Range("M12").Value = 2.650267989565208
$ws.Range("N12").Value = 0.4422666350785391
$ws.Range("B13").Value = 1.078943286621609
$ws.Range("D13").Value = 0.03199453260400276
$ws.Range("E13").Value = 0.06246441301501093
$ws.Range("F13").Value = 0.9760661265253532
$ws.Range("G13").Value = 0.7863779326894758
$ws.Range("H13").Value = 0.118698007071032
$ws.Range("I13").Value = 0.01893508862603888
$ws.Range("J13").Value = 0.5452046883077344
$ws.Range("K13").Value = 0.4815429334279102
$ws.Range("L13").Value = 0.1512264026981782
$ws.Range("M13").Value = 2.683407379581297
$ws.Range("N13").Value = 0.3756329454175642
$ws.Range("B14").Value = 1.040520255483983
$ws.Range("D14").Value = 0.02968556783949339
$ws.Range("E14").Value = 0.06727275095635399
$ws.Range("F14").Value = 0.8793728737668047
$ws.Range("G14").Value = 0.6990997071144562
$ws.Range("H14").Value = 0.1686007439606101
$ws.Range("I14").Value = 0.01999684421118086
$ws.Range("J14").Value = 0.5018592199211014
$ws.Range("K14").Value = 0.4487348753671725
$ws.Range("L14").Value = 0.1709068827960607
$ws.Range("M14").Value = 2.671519739535626
$ws.Range("N14").Value = 0.3272347418247818
$ws.Range("B15").Value = 1.023605651188575
$ws.Range("D15").Value = 0.02901277621600329
$ws.Range("E15").Value = 0.06845488307967251
$ws.Range("F15").Value = 0.8542502799041429
$ws.Range("G15").Value = 0.6760611092293516
$ws.Range("H15").Value = 0.1813857279851021
$ws.Range("I15").Value = 0.02061606044918296
$ws.Range("J15").Value = 0.4909058635509922
$ws.Range("K15").Value = 0.4419508709249325
$ws.Range("L15").Value = 0.1757113274393447
$ws.Range("M15").Value = 2.652819970930722
$ws.Range("N15").Value = 0.3143447902666594
$ws.Range("B16").Value = 0.9602345167040767
$ws.Range("D16").Value = 0.02776330321786702
$ws.Range("E16").Value = 0.06522829798411678
$ws.Range("F16").Value = 0.8458977092249285
$ws.Range("G16").Value = 0.6657776486921705
$ws.Range("H16").Value = 0.1701850881692621
$ws.Range("I16").Value = 0.02300352712438691
$ws.Range("J16").Value = 0.4909578251073867
$ws.Range("K16").Value = 0.4525093081992022
$ws.Range("L16").Value = 0.166393489744415
$ws.Range("M16").Value = 2.485302672334853
$ws.Range("N16").Value = 0.2971702841779944
$ws.Range("B17").Value = 0.9319255797948642
$ws.Range("D17").Value = 0.02800258790216503
$ws.Range("E17").Value = 0.05988710162203148
$ws.Range("F17").Value = 0.8897850204024422
$ws.Range("G17").Value = 0.7033344396474632
$ws.Range("H17").Value = 0.1333279010675881
$ws.Range("I17").Value = 0.02431630337495427
$ws.Range("J17").Value = 0.5136505639748492
$ws.Range("K17").Value = 0.4770252275485198
$ws.Range("L17").Value = 0.1480495984188863
$ws.Range("M17").Value = 2.363843332061634
$ws.Range("N17").Value = 0.309256143703692
$ws.Range("B18").Value = 0.932237270171612
$ws.Range("D18").Value = 0.02966983135745238
$ws.Range("E18").Value = 0.05390755985423823
$ws.Range("F18").Value = 0.9906746932423403
$ws.Range("G18").Value = 0.7927993961180277
$ws.Range("H18").Value = 0.08074281185087528
$ws.Range("I18").Value = 0.02436654457601506
$ws.Range("J18").Value = 0.5615652736851331
$ws.Range("K18").Value = 0.5230104835818281
$ws.Range("L18").Value = 0.1242289889429564
$ws.Range("M18").Value = 2.26495306903567
$ws.Range("N18").Value = 0.348151362700392
$ws.Range("B19").Value = 0.9510421339038828
$ws.Range("D19").Value = 0.03261914913062824
$ws.Range("E19").Value = 0.05072759497470392
$ws.Range("F19").Value = 1.13090644268604
$ws.Range("G19").Value = 0.9178226471200901
$ws.Range("H19").Value = 0.03491473032355685
$ws.Range("I19").Value = 0.02398290852318219
$ws.Range("J19").Value = 0.6261194286188783
$ws.Range("K19").Value = 0.5867808574134408
$ws.Range("L19").Value = 0.1053622075287599
$ws.Range("M19").Value = 2.200338774746712
$ws.Range("N19").Value = 0.411838076558098
$ws.Range("B20").Value = 1.029431427864836
$ws.Range("D20").Value = 0.0390102923934954
$ws.Range("E20").Value = 0.0546107148340198
$ws.Range("F20").Value = 1.387842394446878
$ws.Range("G20").Value = 1.14905841645438
$ws.Range("H20").Value = 0.0061245820339626
$ws.Range("I20").Value = 0.02130312953882196
$ws.Range("J20").Value = 0.7405692892781133
$ws.Range("K20").Value = 0.7034117493320622
$ws.Range("L20").Value = 0.09790706177337771
$ws.Range("M20").Value = 2.218523236826314
$ws.Range("N20").Value = 0.5378496577632887
$ws.Range("B21").Value = 1.160735590504942
$ws.Range("D21").Value = 0.04333533946276802
$ws.Range("E21").Value = 0.05940183298706003
$ws.Range("F21").Value = 1.493018385740029
$ws.Range("G21").Value = 1.248022396071093
$ws.Range("H21").Value = 0.00337090214585789
$ws.Range("I21").Value = 0.01712299288387431
$ws.Range("J21").Value = 0.7808752753878707
$ws.Range("K21").Value = 0.723774918202956
$ws.Range("L21").Value = 0.1052469983679964
$ws.Range("M21").Value = 2.494806598287511
$ws.Range("N21").Value = 0.6136602809446146
$ws.Range("B22").Value = 1.247470349441727
$ws.Range("D22").Value = 0.04590971586068804
$ws.Range("E22").Value = 0.06234693102542943
$ws.Range("F22").Value = 1.556312799919695
$ws.Range("G22").Value = 1.30840916829905
$ws.Range("H22").Value = 0.002126765272699283
$ws.Range("I22").Value = 0.01440275735132346
$ws.Range("J22").Value = 0.8051114909171702
$ws.Range("K22").Value = 0.7351132139048104
$ws.Range("L22").Value = 0.1101327388022
$ws.Range("M22").Value = 2.674821639289206
$ws.Range("N22").Value = 0.6557393153168647
$ws.Range("B23").Value = 1.203997117696673
$ws.Range("D23").Value = 0.04451080886328995
$ws.Range("E23").Value = 0.06086856138434316
$ws.Range("F23").Value = 1.52748937330351
$ws.Range("G23").Value = 1.280943416648341
$ws.Range("H23").Value = 0.002746283870133981
$ws.Range("I23").Value = 0.01544845338111056
$ws.Range("J23").Value = 0.7946007340702579
$ws.Range("K23").Value = 0.7328708183054431
$ws.Range("L23").Value = 0.1075898899408165
$ws.Range("M23").Value = 2.573852885742781
$ws.Range("N23").Value = 0.6308382974379327
$ws.Range("B24").Value = 1.034674632458206
$ws.Range("D24").Value = 0.03933534843223541
$ws.Range("E24").Value = 0.05520879067915718
$ws.Range("F24").Value = 1.412465051171424
$ws.Range("G24").Value = 1.171558440187027
$ws.Range("H24").Value = 0.005893203554333037
$ws.Range("I24").Value = 0.02064689690556953
$ws.Range("J24").Value = 0.7521444933751127
$ws.Range("K24").Value = 0.7185329305159769
$ws.Range("L24").Value = 0.09798400427657583
$ws.Range("M24").Value = 2.201912284410582
$ws.Range("N24").Value = 0.5413514771740608
$ws.Range("B25").Value = 0.8518047782471854
$ws.Range("D25").Value = 0.03389667586302636
$ws.Range("E25").Value = 0.04925280524736309
$ws.Range("F25").Value = 1.293047353799608
$ws.Range("G25").Value = 1.058535599752176
$ws.Range("H25").Value = 0.01082750291146484
$ws.Range("I25").Value = 0.02755540699960157
$ws.Range("J25").Value = 0.7090968124109338
$ws.Range("K25").Value = 0.7039257604889642
$ws.Range("L25").Value = 0.08782308140609807
$ws.Range("M25").Value = 1.805216672005059
$ws.Range("N25").Value = 0.4463945576452772
